# Requirements.xlsx - "environment changed to production"
#
# Updates the Requirements sheet:
#   - Row 2 (Line graph item) gets a Status of "Implemented" and a Comment
#     of "Available in new version."
#   - Row 5 (Investments) description is replaced to reflect that the
#     module still needs to be implemented.
#   - Row 6 (Profile) description is replaced to reflect that the module
#     still needs to be implemented.
#   - The active selection moves to C7.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New Status / Comments for the first requirement row
$ws.Range("D2").Value = "Implemented"
$ws.Range("E2").Value = "Available in new version."

# Updated descriptions for the Investments and Profile rows
$ws.Range("C5").Value = "Investments module need to be implemented."
$ws.Range("C6").Value = "Profile module need to be implemented."

# Move the active selection to C7
$ws.Range("C7").Select()
